$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original export was missing two columns ("Date of Last Update" and a
# duplicate "TicketID") which caused every column from "EH & S Issues"
# onward to be shifted left by one -- the "EH & S Issues" column ended up
# holding a timestamp instead of the real issue category. Insert the two
# missing columns (I and J) before the old "EH & S Issues" column so the
# layout -- and the values that belong in each column -- is corrected.

$ws.Columns("I:J").Insert()

# Headers (row 1)
$ws.Range("I1").Value = "Date of Last Update"
$ws.Range("J1").Value = "TicketID"

# Row 2 (TicketID 241)
$ws.Range("I2").Value = "2017-08-15T09:17:13.210000"
$ws.Range("A2").Copy()
$ws.Range("J2").PasteSpecial(-4163)
$ws.Range("K2").Value = "Other"

# Row 3 (TicketID 238)
$ws.Range("I3").Value = "2017-08-14T16:27:57.897000"
$ws.Range("A3").Copy()
$ws.Range("J3").PasteSpecial(-4163)
$ws.Range("K3").Value = "Other"

$excel.CutCopyMode = 0
